# Add new daily rows (419-449) to each sheet of the Cosenza COVID workbook.
# Rows 419-421 carry new case/death/discharge/hospitalisation/ICU counts plus
# the rolling 7-day AVERAGE() formulas (and, on "Ricoveri", the day-over-day
# delta formulas). Rows 422-449 only carry the date in column A (data not
# yet available for those future days at the time of the commit).

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# New C-column values per sheet (by sheet name) for rows 419, 420, 421.
$newValues = @{
    "Nuovi casi"        = @(167, 155, 91)
    "Deceduti"           = @(3, 2, 6)
    "Dimessi   Guariti"  = @(295, 0, 723)
    "Ricoveri"           = @(184, 180, 179)
    "Terapia"            = @(17, 18, 17)
}

$sheetNames = @("Nuovi casi", "Deceduti", "Dimessi   Guariti", "Ricoveri", "Terapia")

$firstNewRow = 419
$lastDataRow = 421
$lastDateRow = 449
$firstDate = 44317

# Column-A date cells in every sheet's new rows use the "dd/mm/yyyy, default
# font" style (the one already used for column A throughout "Nuovi casi",
# "Deceduti" and "Dimessi   Guariti") -- even on "Ricoveri"/"Terapia", whose
# *older* rows use a slightly different (explicit black font) date style.
# Grab a reference cell with that exact style to stamp onto every new A cell.
$dateStyleSource = $wb.Worksheets.Item("Nuovi casi").Range("A418")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $values = $newValues[$name]

    for ($row = $firstNewRow; $row -le $lastDateRow; $row++) {

        # --- column A: date, using the shared "A419:A449" date style ---
        $dateStyleSource.Copy() | Out-Null
        $ws.Range("A" + $row).PasteSpecial($xlPasteFormats) | Out-Null
        $ws.Range("A" + $row).Value = $firstDate + ($row - $firstNewRow)

        if ($row -le $lastDataRow) {
            $idx = $row - $firstNewRow

            # --- column C: reported count ---
            $ws.Range("C" + ($row - 1)).Copy() | Out-Null
            $ws.Range("C" + $row).PasteSpecial($xlPasteFormats) | Out-Null
            $ws.Range("C" + $row).Value = $values[$idx]

            # --- column D: rolling 7-day average, same format as row above ---
            $ws.Range("D" + ($row - 1)).Copy() | Out-Null
            $ws.Range("D" + $row).PasteSpecial($xlPasteFormats) | Out-Null
            $ws.Range("D" + $row).Formula = "=AVERAGE(C" + ($row - 6) + ":C" + $row + ")"

            # --- column E (Ricoveri only): day-over-day delta ---
            if ($name -eq "Ricoveri") {
                $ws.Range("E" + ($row - 1)).Copy() | Out-Null
                $ws.Range("E" + $row).PasteSpecial($xlPasteFormats) | Out-Null
                $ws.Range("E" + $row).Formula = "=C" + $row + "-C" + ($row - 1)
            }
        }
    }

    $excel.CutCopyMode = 0
}

# --- restore view / selection state seen in the saved workbook ---
$wsNuoviCasi = $wb.Worksheets.Item("Nuovi casi")
$wsNuoviCasi.Activate()
$excel.ActiveWindow.ScrollRow = 403
$wsNuoviCasi.Range("C419:C421").Select() | Out-Null

$wsDeceduti = $wb.Worksheets.Item("Deceduti")
$wsDeceduti.Activate()
$excel.ActiveWindow.ScrollRow = 405
$wsDeceduti.Range("C419:C421").Select() | Out-Null

$wsDimessi = $wb.Worksheets.Item("Dimessi   Guariti")
$wsDimessi.Activate()
$excel.ActiveWindow.ScrollRow = 404
$wsDimessi.Range("C419:C421").Select() | Out-Null

$wsRicoveri = $wb.Worksheets.Item("Ricoveri")
$wsRicoveri.Activate()
$excel.ActiveWindow.ScrollRow = 404
$wsRicoveri.Range("C419:C421").Select() | Out-Null

$wsTerapia = $wb.Worksheets.Item("Terapia")
$wsTerapia.Activate()
$excel.ActiveWindow.ScrollRow = 396
$wsTerapia.Range("F413").Select() | Out-Null
